$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44315
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("S2").Value = 1333
$ws.Range("D3").Value = 44315
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 1111
$ws.Range("D4").Value = 44279
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 50
$ws.Range("R4").Value = 'Provincia de Melipilla'
$ws.Range("D5").Value = 44279
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Provincia de Melipilla'
$ws.Range("S5").Value = 667
$ws.Range("T5").Value = 18
$ws.Range("D6").Value = 44277
$ws.Range("L6").Value = 'Especial'
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 833
$ws.Range("D7").Value = 44291
$ws.Range("L7").Value = 'Extra (doble especial)'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("S7").Value = 1000
$ws.Range("D8").Value = 44350
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 24000
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1333
$ws.Range("D9").Value = 44273
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = '$/caja 16 kilos'
$ws.Range("R9").Value = 'Provincia de Melipilla'
$ws.Range("S9").Value = 938
$ws.Range("T9").Value = 16
$ws.Range("D10").Value = 44273
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 13000
$ws.Range("Q10").Value = '$/caja 16 kilos'
$ws.Range("R10").Value = 'Provincia de Melipilla'
$ws.Range("S10").Value = 812
$ws.Range("T10").Value = 16
$ws.Range("D11").Value = 44273
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("Q11").Value = '$/caja 16 kilos'
$ws.Range("S11").Value = 625
$ws.Range("T11").Value = 16
$ws.Range("D12").Value = 44300
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 120
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 16000
$ws.Range("S13").Value = 889
$ws.Range("D14").Value = 44222
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("Q14").Value = '$/caja 16 kilos'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 1125
$ws.Range("T14").Value = 16
$ws.Range("D15").Value = 44630
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = '$/caja 20 kilos'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 20
$ws.Range("D16").Value = 44298
$ws.Range("L16").Value = 'Extra (doble especial)'
$ws.Range("M16").Value = 160
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 20000
$ws.Range("S16").Value = 1111
$ws.Range("D17").Value = 44267
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 13000
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 13000
$ws.Range("S17").Value = 722
$ws.Range("D18").Value = 44295
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 130
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Provincia de Melipilla'
$ws.Range("S18").Value = 556
$ws.Range("T18").Value = 18
$ws.Range("D19").Value = 44292
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 150
$ws.Range("D20").Value = 44292
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 80
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("S20").Value = 778
$ws.Range("T20").Value = 18
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 14000
$ws.Range("O21").Value = 14000
$ws.Range("P21").Value = 14000
$ws.Range("S21").Value = 875
$ws.Range("D22").Value = 44274
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 130
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("Q22").Value = '$/caja 16 kilos'
$ws.Range("S22").Value = 750
$ws.Range("T22").Value = 16
$ws.Range("D23").Value = 44645
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 16000
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 16000
$ws.Range("R23").Value = 'Provincia de Limarí'
$ws.Range("S23").Value = 889
$ws.Range("D24").Value = 44271
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 15000
$ws.Range("P24").Value = 15000
$ws.Range("S24").Value = 833
$ws.Range("D25").Value = 44258
$ws.Range("L25").Value = 'Primera'
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("S25").Value = 778
$ws.Range("D26").Value = 44299
$ws.Range("L26").Value = 'Especial'
$ws.Range("M26").Value = 170
$ws.Range("N26").Value = 18000
$ws.Range("O26").Value = 18000
$ws.Range("P26").Value = 18000
$ws.Range("S26").Value = 1000
$ws.Range("D27").Value = 44299
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 16000
$ws.Range("O27").Value = 16000
$ws.Range("P27").Value = 16000
$ws.Range("S27").Value = 889
$ws.Range("D28").Value = 44284
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 13000
$ws.Range("O28").Value = 13000
$ws.Range("P28").Value = 13000
$ws.Range("S28").Value = 722
$ws.Range("D29").Value = 44284
$ws.Range("L29").Value = 'Extra (doble especial)'
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 15000
$ws.Range("S29").Value = 833
$ws.Range("D30").Value = 44284
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 12000
$ws.Range("O30").Value = 12000
$ws.Range("P30").Value = 12000
$ws.Range("Q30").Value = '$/caja 18 kilos'
$ws.Range("S30").Value = 667
$ws.Range("T30").Value = 18
$ws.Range("D31").Value = 44301
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = '$/caja 18 kilos'
$ws.Range("S31").Value = 889
$ws.Range("T31").Value = 18
$ws.Range("D32").Value = 44224
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 120
$ws.Range("N32").Value = 18000
$ws.Range("O32").Value = 18000
$ws.Range("P32").Value = 18000
$ws.Range("R32").Value = 'Provincia de Limarí'
$ws.Range("S32").Value = 1125
$ws.Range("D33").Value = 44252
$ws.Range("M33").Value = 140
$ws.Range("N33").Value = 13000
$ws.Range("O33").Value = 13000
$ws.Range("P33").Value = 13000
$ws.Range("S33").Value = 722
$ws.Range("D34").Value = 44309
$ws.Range("L34").Value = 'Especial'
$ws.Range("M34").Value = 100
$ws.Range("D35").Value = 44309
$ws.Range("M35").Value = 60
$ws.Range("N35").Value = 18000
$ws.Range("O35").Value = 18000
$ws.Range("P35").Value = 18000
$ws.Range("S35").Value = 1000
$ws.Range("D36").Value = 44330
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 50
$ws.Range("N36").Value = 23000
$ws.Range("O36").Value = 23000
$ws.Range("P36").Value = 23000
$ws.Range("R36").Value = 'Provincia de Melipilla'
$ws.Range("S36").Value = 1278
